# RKI-COVID-19_Todesfaelle.xlsx update (download 2020-12-15--13-35-01)
# Sheet "COVID_Todesfälle"      (weekly death counts, column B is TEXT)
# Sheet "COVID_Todesfälle_Monat" (monthly death counts, column B is TEXT)
#
# Column B on both sheets holds numeric-looking values stored as TEXT
# (shared strings) in the workbook, so every write below forces the
# cell to Text format first (so Excel doesn't auto-convert "1740" etc.
# into a number) and then clears the number-format back off the cell
# once the text value has been committed, leaving the cell on the
# default/general style exactly like its untouched neighbours.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---------------------------------------------------------------
# Sheet 1: COVID_Todesfälle
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("COVID_Todesfälle")

$s1Updates = @{
    5  = "<4"
    16 = "1740"
    19 = "783"
    20 = "514"
    21 = "351"
    22 = "271"
    24 = "112"
    25 = "73"
    34 = "32"
    38 = "32"
    39 = "52"
    42 = "115"
    43 = "226"
    44 = "378"
    45 = "725"
}

foreach ($row in $s1Updates.Keys) {
    Set-TextValue $ws1.Cells.Item($row, 2) $s1Updates[$row]
}

# Two new weeks appended at the bottom of the table.
$s1NewRows = @(
    @{ Row = 46; A = 45; B = "1105" },
    @{ Row = 47; A = 46; B = "1455" }
)

foreach ($r in $s1NewRows) {
    $ws1.Cells.Item($r.Row, 1).Value = $r.A
    Set-TextValue $ws1.Cells.Item($r.Row, 2) $r.B
}

# ---------------------------------------------------------------
# Sheet 2: COVID_Todesfälle_Monat
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("COVID_Todesfälle_Monat")

$s2Updates = @{
    2 = "1117"
    3 = "6041"
    4 = "1562"
    7 = "145"
    8 = "200"
    9 = "1366"
}

foreach ($row in $s2Updates.Keys) {
    Set-TextValue $ws2.Cells.Item($row, 2) $s2Updates[$row]
}

# One new month appended at the bottom of the table.
$ws2.Cells.Item(10, 1).Value = 11
Set-TextValue $ws2.Cells.Item(10, 2) "2682"
